$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 1 (the old helper/numeric cell that used the bold+bordered style),
# which shifts row 2 (holding the JSON text) up into row 1.
$ws.Rows("1:1").Delete()

# Replace the cell text with the reformatted JSON-style text.
$ws.Range("A1").Value = 'questions = [
    {
        "title": "You have a customer on the phone who you are struggling to understand due to their limited proficiency in your language. You need to offer a troubleshooting step for a device issue.These four actions will all have the same outcome:Pressing the power button twice.Shutting down the device and running a backup.Powering off the device and turning it back on.Restarting the system.What should you say to the customer?",
        "ques_type": 2,
        "options": [
            "\u201cPress the big power button two times.\u201d",
            "\u201cShut the device down, then perform a backup.\u201d",
            "\u201cTurn the device off, wait a moment, then turn it back on.\u201d",
            "\u201cClick on restart in the power menu.\u201d"
        ],
        "score": "\u201cPress the big power button two times.\u201d"
    },
    {
        "title": "You are a call center representative for a bank. You receive the following call:Once you have gathered all the necessary details, how should you confirm that the account has been frozen?",
        "ques_type": 2,
        "options": [
            "\u201cI\u2019ve now frozen your account. Thank you for letting us know about this. Is there anything else I can help you with today?\u201d",
            "\u201cI\u2019ve now frozen your account. By the way, have you called the police? What did they say? Did you identify your aggressor?\u201d",
            "\u201cI\u2019ve now frozen your account. I&#39m sorry that you had this experience. Don&#39t hesitate to call if you need any other help with your account.\u201d",
            "\u201cHow awful! I can\u2019t believe that there are people out there who would act in this way. I\u2019ve now frozen your account.\u201d"
        ],
        "score": "\u201cI\u2019ve now frozen your account. I&#39m sorry that you had this experience. Don&#39t hesitate to call if you need any other help with your account.\u201d"
    },
    {
        "title": "You are a call center representative for a technology company. Some of your products are sold with a free additional power cable as a surprise gift. You have been provided the flow chart shown below.You receive the following call:Unfortunately, you do not have direct access to the gift allocation information.\u00a0How should you reply?",
        "ques_type": 2,
        "options": [
            "\u201cOK, I\u2019ll be happy to arrange delivery of a replacement cable for you immediately.\u201d",
            "\u201cLet me check whether an extra power cable was included in your order, then I\u2019ll call you back.\u201d",
            "\u201cCould you try using a different power outlet to see if that solves the problem?\u201d",
            "\u201cCould you check the delivery box to see if an extra power cable was included as a free gift?\u201d"
        ],
        "score": "\u201cCould you check the delivery box to see if an extra power cable was included as a free gift?\u201d"
    },
    {
        "title": "You are a call center representative at a mid-sized IT services company. A customer on the line has an issue with a missed delivery, and expects a quick resolution. However, before you have the time to get the required details from them, they complain to you about your colleague Jo, citing unresolved issues and poor service during a call last week. You''ve noticed a pattern of complaints about Jo, but you''re not the team manager.Which action should you take?",
        "ques_type": 2,
        "options": [
            "Resolve the customer\u2019s delivery issue, and take no further action.",
            "Resolve the customer\u2019s delivery issue, then report the recurring complaints about Jo to your manager.",
            "Apologize to the customer and suggest they submit their feedback to Jo directly, then resolve their delivery issue.",
            "Resolve the customer\u2019s delivery issue, then escalate the feedback about Jo directly to senior management."
        ],
        "score": "Resolve the customer\u2019s delivery issue, then report the recurring complaints about Jo to your manager."
    }
]'

# Undo the automatic row-height expansion triggered by the long text so the
# row keeps the sheet's default height (matches the target workbook).
$ws.Rows(1).AutoFit()
